# Auto-generated edit script: updates LeveProfit market-price columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1014.1667
$ws.Range("J28").Value = 932.4
$ws.Range("L28").Value = 932.4
$ws.Range("N28").Value = -1902.4

$ws.Range("H40").Value = 4155.8823
$ws.Range("J40").Value = 5400.1113
$ws.Range("L40").Value = 5400.1113
$ws.Range("N40").Value = -5750.1113

$ws.Range("H97").Value = 2881.875
$ws.Range("J97").Value = 2881.875
$ws.Range("L97").Value = 8645.625
$ws.Range("N97").Value = -9637.625

$ws.Range("H98").Value = 1488.7
$ws.Range("J98").Value = 2359.4
$ws.Range("L98").Value = 2359.4
$ws.Range("N98").Value = -5355.4

$ws.Range("H107").Value = 516.75
$ws.Range("I107").Value = 404.42856
$ws.Range("J107").Value = 1303
$ws.Range("K107").Value = 404.42856
$ws.Range("L107").Value = 1303
$ws.Range("M107").Value = 1515.57144
$ws.Range("N107").Value = -5143

$ws.Range("H115").Value = 554.3333
$ws.Range("J115").Value = 1000
$ws.Range("L115").Value = 3000
$ws.Range("N115").Value = -6134

$ws.Range("H116").Value = 4307.5835
$ws.Range("I116").Value = 4069.2
$ws.Range("K116").Value = 4069.2
$ws.Range("M116").Value = -627.1999999999998

$ws.Range("H122").Value = 1488.7
$ws.Range("J122").Value = 2359.4
$ws.Range("L122").Value = 7078.200000000001
$ws.Range("N122").Value = -11978.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2162
$ws.Range("I2").Value = 990.7143
$ws.Range("K2").Value = 990.7143
$ws.Range("M2").Value = -877.7143

$ws.Range("H97").Value = 1815.4375
$ws.Range("I97").Value = 1234.4615
$ws.Range("K97").Value = 1234.4615
$ws.Range("M97").Value = -738.4614999999999

$ws.Range("H102").Value = 933.8
$ws.Range("J102").Value = 800
$ws.Range("L102").Value = 800
$ws.Range("N102").Value = -4044

$ws.Range("H110").Value = 538.55554
$ws.Range("I110").Value = 538.55554
$ws.Range("K110").Value = 538.55554
$ws.Range("M110").Value = 1506.44446

$ws.Range("H116").Value = 2162
$ws.Range("I116").Value = 990.7143
$ws.Range("K116").Value = 990.7143
$ws.Range("M116").Value = 1303.2857

$ws.Range("H132").Value = 3429.5715
$ws.Range("I132").Value = 3429.5715
$ws.Range("K132").Value = 10288.7145
$ws.Range("M132").Value = -7758.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2162
$ws.Range("I3").Value = 990.7143
$ws.Range("K3").Value = 990.7143
$ws.Range("M3").Value = -876.7143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3104.4
$ws.Range("I16").Value = 761
$ws.Range("J16").Value = 4666.6665
$ws.Range("K16").Value = 761
$ws.Range("L16").Value = 4666.6665
$ws.Range("M16").Value = -474
$ws.Range("N16").Value = -5240.6665

$ws.Range("H22").Value = 3994.8823
$ws.Range("I22").Value = 4736.643
$ws.Range("J22").Value = 533.3333
$ws.Range("K22").Value = 4736.643
$ws.Range("L22").Value = 533.3333
$ws.Range("M22").Value = -4386.643
$ws.Range("N22").Value = -1233.3333

$ws.Range("H58").Value = 5553.8335
$ws.Range("I58").Value = 5553
$ws.Range("J58").Value = 5554.25
$ws.Range("K58").Value = 5553
$ws.Range("L58").Value = 5554.25
$ws.Range("M58").Value = -5350
$ws.Range("N58").Value = -5960.25

$ws.Range("H99").Value = 2624.2273
$ws.Range("I99").Value = 2691.65
$ws.Range("J99").Value = 1950
$ws.Range("K99").Value = 2691.65
$ws.Range("L99").Value = 1950
$ws.Range("M99").Value = -1193.65
$ws.Range("N99").Value = -4946

$ws.Range("H107").Value = 1370.2333
$ws.Range("I107").Value = 1127.6
$ws.Range("K107").Value = 1127.6
$ws.Range("M107").Value = 792.4000000000001

$ws.Range("H113").Value = 3104.4
$ws.Range("I113").Value = 761
$ws.Range("J113").Value = 4666.6665
$ws.Range("K113").Value = 761
$ws.Range("L113").Value = 4666.6665
$ws.Range("M113").Value = 1409
$ws.Range("N113").Value = -9006.666499999999

$ws.Range("H126").Value = 2624.2273
$ws.Range("I126").Value = 2691.65
$ws.Range("J126").Value = 1950
$ws.Range("K126").Value = 8074.950000000001
$ws.Range("L126").Value = 5850
$ws.Range("M126").Value = -5604.950000000001
$ws.Range("N126").Value = -10790

$ws.Range("H132").Value = 2353.5293
$ws.Range("I132").Value = 2267.3333
$ws.Range("K132").Value = 6801.999899999999
$ws.Range("M132").Value = -4271.999899999999

$ws.Range("H134").Value = 2144
$ws.Range("I134").Value = 2189.3
$ws.Range("K134").Value = 6567.900000000001
$ws.Range("M134").Value = -4032.900000000001

$ws.Range("H136").Value = 5553.8335
$ws.Range("I136").Value = 5553
$ws.Range("J136").Value = 5554.25
$ws.Range("K136").Value = 16659
$ws.Range("L136").Value = 16662.75
$ws.Range("M136").Value = -14109
$ws.Range("N136").Value = -21762.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H121").Value = 196.66667
$ws.Range("I121").Value = 196.66667
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 590.00001
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = 719.99999
$ws.Range("M121").ClearContents()

$ws.Range("H131").Value = 1500
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1500
$ws.Range("K131").Value = 0
$ws.Range("M131").Value = 4500
$ws.Range("N131").Value = -14580
$ws.Range("L131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 1625
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 1625
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = 1625
$ws.Range("N29").Value = -2205
$ws.Range("L29").ClearContents()

$ws.Range("H107").Value = 4981.727
$ws.Range("I107").Value = 433.33334
$ws.Range("J107").Value = 10439.8
$ws.Range("K107").Value = 433.33334
$ws.Range("L107").Value = 10439.8
$ws.Range("M107").Value = 1486.66666
$ws.Range("N107").Value = -14279.8

$ws.Range("H113").Value = 2699.3333
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 4583
$ws.Range("I132").Value = 5624.5
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 16873.5
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -14343.5
$ws.Range("N132").Value = -12560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2733
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 3599.5
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 3599.5
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -4189.5

$ws.Range("H27").Value = 2733
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 3599.5
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 3599.5
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -3813.5

$ws.Range("H46").Value = 3813.3784
$ws.Range("I46").Value = 3396.1482
$ws.Range("J46").Value = 4939.9
$ws.Range("K46").Value = 3396.1482
$ws.Range("L46").Value = 4939.9
$ws.Range("M46").Value = -3208.1482
$ws.Range("N46").Value = -5315.9

$ws.Range("H61").Value = 2361.3
$ws.Range("I61").Value = 2487.5715
$ws.Range("J61").Value = 2066.6667
$ws.Range("K61").Value = 2487.5715
$ws.Range("L61").Value = 2066.6667
$ws.Range("M61").Value = -2285.5715
$ws.Range("N61").Value = -2470.6667

$ws.Range("H113").Value = 2361.3
$ws.Range("I113").Value = 2487.5715
$ws.Range("J113").Value = 2066.6667
$ws.Range("K113").Value = 2487.5715
$ws.Range("L113").Value = 2066.6667
$ws.Range("M113").Value = -317.5715
$ws.Range("N113").Value = -6406.6667

$ws.Range("H132").Value = 4780.875
$ws.Range("I132").Value = 4672
$ws.Range("J132").Value = 4962.3335
$ws.Range("K132").Value = 14016
$ws.Range("L132").Value = 14887.0005
$ws.Range("M132").Value = -11486
$ws.Range("N132").Value = -19947.0005

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 264
$ws.Range("I107").Value = 238.72728
$ws.Range("K107").Value = 716.18184
$ws.Range("M107").Value = 1203.81816

$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -470
$ws.Range("M132").ClearContents()

$ws.Range("H136").Value = 1014.4
$ws.Range("I136").Value = 976.4167
$ws.Range("K136").Value = 2929.2501
$ws.Range("M136").Value = -379.2501000000002
